$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (A1:E1): shade with light-gray fill and add a thin top/bottom
#    rule (keeps the existing bold, center/center Times New Roman font).
# ---------------------------------------------------------------------------
$hdr = $ws.Range("A1:E1")

$hdr.Interior.Color = 14277081   # RGB(217,217,217) = D9D9D9

$hdr.Borders.Item(8).LineStyle = 1    # xlEdgeTop, xlContinuous (thin)
$hdr.Borders.Item(8).Color = 0
$hdr.Borders.Item(9).LineStyle = 1    # xlEdgeBottom, xlContinuous (thin)
$hdr.Borders.Item(9).Color = 0

# ---------------------------------------------------------------------------
# 2. New footnote row (row 9), merged A9:E9, 8pt italic Times New Roman,
#    left/center + wrap, double top/bottom rule, 45pt row height.
# ---------------------------------------------------------------------------
$footnote = "*All values displayed as mean ± SD for ratio continuous variables or n (%) for dichotomous categorical variables. Percentages for the variant columns were calculated in respect to total patients within a variant (i.e., within column), and percentages for the total column was calculated in respect to the population total."

$note = $ws.Range("A9:E9")

# Start from the existing italic body-label style (A5) so we only have to
# change the font size and add wrap/border on top of it.
$ws.Range("A5").Copy()
$note.PasteSpecial(-4122)   # xlPasteFormats

$note.Value = $footnote
$note.Font.Size = 8
$note.WrapText = $true

$note.Borders.Item(8).LineStyle = -4119   # xlEdgeTop, xlDouble
$note.Borders.Item(8).Color = 0
$note.Borders.Item(9).LineStyle = -4119   # xlEdgeBottom, xlDouble
$note.Borders.Item(9).Color = 0

$note.Merge()
$ws.Rows.Item(9).RowHeight = 45

Write-Host "done"
